$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure D-column price cells are treated as text (matches source formatting, avoids numeric coercion)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.837.37"
$ws.Range("E2").Value = "  +1.09%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.623.54"
$ws.Range("E3").Value = "  +0.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  -0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.96"
$ws.Range("E5").Value = "  -1.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "151.58"
$ws.Range("E6").Value = "  +3.98%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.00%  "

$ws.Range("E8").Value = "  +0.93%  "

$ws.Range("E9").Value = "  +1.91%  "

$ws.Range("E10").Value = "  +3.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.386"
$ws.Range("E11").Value = "  +4.17%  "

$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.91"
$ws.Range("E13").Value = "  +2.94%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.090.87"
$ws.Range("E14").Value = "  +0.46%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "63.647.20"
$ws.Range("E15").Value = "  +1.12%  "

$ws.Range("E16").Value = "  +5.04%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "2.607.94"
$ws.Range("E17").Value = "  -0.58%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.40"
$ws.Range("E18").Value = "  +8.28%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.73"
$ws.Range("E19").Value = "  +4.72%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "349.11"
$ws.Range("E20").Value = "  +2.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  +1.04%  "

$ws.Range("E22").Value = "  -0.34%  "

# Rows 23-51: coin rows shift up by one; update B (name), C (link), D (price), E (change%)
$ws.Range("B23").Value = "Litecoin"
$ws.Range("C23").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.12"
$ws.Range("E23").Value = "  +1.50%  "

$ws.Range("B24").Value = "SuiNetwork"
$ws.Range("C24").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.73"
$ws.Range("E24").Value = "  +8.78%  "

$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.41"
$ws.Range("E25").Value = "  +4.47%  "

$ws.Range("B26").Value = "Fetch.AI"
$ws.Range("C26").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.68"
$ws.Range("E26").Value = "  -0.39%  "

$ws.Range("B27").Value = "Bittensor"
$ws.Range("C27").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "558.29"
$ws.Range("E27").Value = "  -0.19%  "

$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.11"
$ws.Range("E28").Value = "  +4.46%  "

$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.163"
$ws.Range("E29").Value = "  +0.91%  "

$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.04%  "

$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.07"
$ws.Range("E31").Value = "  +2.07%  "

$ws.Range("B32").Value = "PEPE"
$ws.Range("C32").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0853"
$ws.Range("E32").Value = "  +1.20%  "

$ws.Range("B33").Value = "ImmutableX"
$ws.Range("C33").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.75"
$ws.Range("E33").Value = "  +0.30%  "

$ws.Range("B34").Value = "NEARProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.32"
$ws.Range("E34").Value = "  +3.15%  "

$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "167.53"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.418"
$ws.Range("E36").Value = "  +4.22%  "

$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.10%  "

$ws.Range("B38").Value = "EthereumClassic"
$ws.Range("C38").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.63"
$ws.Range("E38").Value = "  +3.49%  "

$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.95"
$ws.Range("E39").Value = "  +1.56%  "

$ws.Range("B40").Value = "USDe"
$ws.Range("C40").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.999"
$ws.Range("E40").Value = "  -0.02%  "

$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "168.17"
$ws.Range("E41").Value = "  +2.18%  "

$ws.Range("B42").Value = "OKB"
$ws.Range("C42").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "39.78"
$ws.Range("E42").Value = "  +0.50%  "

$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.98"
$ws.Range("E43").Value = "  +6.40%  "

$ws.Range("B44").Value = "Hedera"
$ws.Range("C44").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0591"
$ws.Range("E44").Value = "  +4.82%  "

$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "21.95"
$ws.Range("E45").Value = "  +1.06%  "

$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.633"
$ws.Range("E46").Value = "  +1.39%  "

$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0253"
$ws.Range("E47").Value = "  +3.23%  "

$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.04"
$ws.Range("E48").Value = "  +6.54%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0969"
$ws.Range("E49").Value = "  +1.42%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "19.42"
$ws.Range("E50").Value = "  +4.25%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0241"
$ws.Range("E51").Value = "  +22.40%  "
